# Update 'F' column ( 想去人数 / want-to-go count ) values across all 4 sheets
# per commit: "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 425  # F2: 418 -> 425
$ws.Cells.Item(3, 6).Value = 579  # F3: 573 -> 579
$ws.Cells.Item(4, 6).Value = 880  # F4: 878 -> 880
$ws.Cells.Item(5, 6).Value = 624  # F5: 614 -> 624
$ws.Cells.Item(6, 6).Value = 793  # F6: 790 -> 793
$ws.Cells.Item(7, 6).Value = 362  # F7: 360 -> 362
$ws.Cells.Item(8, 6).Value = 552  # F8: 551 -> 552
$ws.Cells.Item(9, 6).Value = 110  # F9: 109 -> 110
$ws.Cells.Item(10, 6).Value = 1122  # F10: 1121 -> 1122
$ws.Cells.Item(11, 6).Value = 571  # F11: 570 -> 571
$ws.Cells.Item(12, 6).Value = 334  # F12: 333 -> 334
$ws.Cells.Item(13, 6).Value = 451  # F13: 448 -> 451
$ws.Cells.Item(15, 6).Value = 291  # F15: 289 -> 291
$ws.Cells.Item(16, 6).Value = 43  # F16: 41 -> 43
$ws.Cells.Item(18, 6).Value = 524  # F18: 523 -> 524
$ws.Cells.Item(19, 6).Value = 14  # F19: 7 -> 14
$ws.Cells.Item(20, 6).Value = 520  # F20: 516 -> 520
$ws.Cells.Item(21, 6).Value = 10  # F21: 9 -> 10
$ws.Cells.Item(22, 6).Value = 470  # F22: 465 -> 470

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 67  # F2: 66 -> 67
$ws.Cells.Item(3, 6).Value = 53  # F3: 52 -> 53
$ws.Cells.Item(4, 6).Value = 304  # F4: 303 -> 304
$ws.Cells.Item(7, 6).Value = 635  # F7: 634 -> 635
$ws.Cells.Item(10, 6).Value = 43  # F10: 42 -> 43
$ws.Cells.Item(12, 6).Value = 19  # F12: 18 -> 19

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 363  # F2: 361 -> 363

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 363  # F2: 361 -> 363
$ws.Cells.Item(3, 6).Value = 67  # F3: 66 -> 67
$ws.Cells.Item(4, 6).Value = 426  # F4: 418 -> 426
$ws.Cells.Item(5, 6).Value = 53  # F5: 52 -> 53
$ws.Cells.Item(6, 6).Value = 304  # F6: 303 -> 304
$ws.Cells.Item(7, 6).Value = 579  # F7: 573 -> 579
$ws.Cells.Item(8, 6).Value = 880  # F8: 878 -> 880
$ws.Cells.Item(9, 6).Value = 624  # F9: 614 -> 624
$ws.Cells.Item(10, 6).Value = 793  # F10: 790 -> 793
$ws.Cells.Item(11, 6).Value = 362  # F11: 360 -> 362
$ws.Cells.Item(12, 6).Value = 552  # F12: 551 -> 552
$ws.Cells.Item(13, 6).Value = 110  # F13: 109 -> 110
$ws.Cells.Item(14, 6).Value = 1122  # F14: 1121 -> 1122
$ws.Cells.Item(15, 6).Value = 571  # F15: 570 -> 571
$ws.Cells.Item(18, 6).Value = 334  # F18: 333 -> 334
$ws.Cells.Item(19, 6).Value = 451  # F19: 448 -> 451
$ws.Cells.Item(20, 6).Value = 635  # F20: 634 -> 635
$ws.Cells.Item(23, 6).Value = 291  # F23: 289 -> 291
$ws.Cells.Item(24, 6).Value = 43  # F24: 41 -> 43
$ws.Cells.Item(27, 6).Value = 43  # F27: 42 -> 43
$ws.Cells.Item(28, 6).Value = 524  # F28: 523 -> 524
$ws.Cells.Item(30, 6).Value = 19  # F30: 18 -> 19
$ws.Cells.Item(31, 6).Value = 14  # F31: 7 -> 14
$ws.Cells.Item(32, 6).Value = 520  # F32: 516 -> 520
$ws.Cells.Item(33, 6).Value = 10  # F33: 9 -> 10
$ws.Cells.Item(34, 6).Value = 470  # F34: 465 -> 470
